$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score Cards")  # same sheet as $wb.ActiveSheet

# New round played on 28 June 2025 at Ocean View.
$ws.Range("A12").Value = "Ocean View"
$ws.Range("B12").Value = 45836
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)

$holeScores = @(5, 4, 6, 4, 6, 4, 6, 4, 4, 5, 4, 8, 7, 5, 5, 5, 4, 4)
for ($i = 0; $i -lt $holeScores.Length; $i++) {
    $col = 3 + $i  # Column C = hole 1
    $ws.Cells.Item(12, $col).Value = $holeScores[$i]
}

# GIR, Putts, Fairways
$ws.Range("U12").Value = 3
$ws.Range("V12").Value = 42
$ws.Range("W12").Value = 5

# Grow Table1 to include the newly entered row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:W12"))

# Move the selection to mirror where Excel leaves the cursor afterward.
$ws.Range("V13").Select()
